$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-03 Friday" "2024-05-04 Saturday"

Replace-Text "106×9=954" "614×8=4912"
Replace-Text "848×9=7632" "359×8=2872"
Replace-Text "901×7=6307" "848×3=2544"
Replace-Text "955×8=7640" "750×7=5250"
Replace-Text "515×9=4635" "766×9=6894"

Replace-Text "734×7=5138" "344×6=2064"
Replace-Text "872×4=3488" "991×2=1982"
Replace-Text "428×3=1284" "837×2=1674"
Replace-Text "505×6=3030" "294×5=1470"
Replace-Text "428×2=856" "904×5=4520"

Replace-Text "423×8=3384" "255×9=2295"
Replace-Text "985×8=7880" "151×6=906"
Replace-Text "480×3=1440" "921×8=7368"
Replace-Text "286×2=572" "953×2=1906"
Replace-Text "876×3=2628" "778×2=1556"

Replace-Text "571×4=2284" "137×9=1233"
Replace-Text "534×4=2136" "473×3=1419"
Replace-Text "539×7=3773" "166×3=498"
Replace-Text "993×2=1986" "849×4=3396"
Replace-Text "465×4=1860" "506×2=1012"

Replace-Text "418×7=2926" "266×5=1330"
Replace-Text "987×6=5922" "124×8=992"
Replace-Text "296×8=2368" "166×4=664"
Replace-Text "103×3=309" "907×9=8163"
Replace-Text "909×6=5454" "309×8=2472"
